{"js": "// The presentation rubric refers to reference formatting as \"APA 7th\n// addition\" in four places; this should read \"APA 7th edition\". The\n// last occurrence (in the \"does not fully contain evidence\" column)\n// is where the document's \"_GoBack\" bookmark (Word's marker for the\n// most recent edit location) lives, immediately after the word and\n// before the trailing period, so it is moved there too.\n\n// 1) Drop the existing \"_GoBack\" bookmark - it will be re-inserted at\n//    its new location once we know where the last replacement landed.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) Replace every whole-word, case-sensitive occurrence of \"addition\"\n//    with \"edition\" (there are exactly four, one per rubric column).\nconst body = context.document.body;\nconst additionMatches = body.search(\"addition\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nadditionMatches.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < additionMatches.items.length; i++) {\n  additionMatches.items[i].insertText(\"edition\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) Re-find \"edition\" and re-insert the bookmark as a collapsed range\n//    right after the fourth (last) occurrence, i.e. before the \".\".\nconst editionMatches = body.search(\"edition\", {\n  matchCase: true,\n  matchWholeWord: true\n});\neditionMatches.load(\"text\");\nawait context.sync();\n\nconst lastEdition = editionMatches.items[editionMatches.items.length - 1];\nconst insertionPoint = lastEdition.getRange(Word.RangeLocation.end);\ninsertionPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The presentation rubric refers to reference formatting as \"APA 7th\n# addition\" in four places; this should read \"APA 7th edition\". The\n# last occurrence (in the \"does not fully contain evidence\" column) is\n# where the document's \"_GoBack\" bookmark (Word's marker for the most\n# recent edit location) lives, immediately after the word and before\n# the trailing period, so it is moved there too.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the existing \"_GoBack\" bookmark - it will be re-added at its\n#    new location once we know where the last replacement landed.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Replace every whole-word, case-sensitive occurrence of \"addition\"\n#    with \"edition\" (there are exactly four, one per rubric column).\n$find = $d.Content.Find\n$find.Execute(\"addition\", $true, $true, $false, $false, $false, $true, 1, $false, \"edition\", 2)\n\n# 3) Re-find \"edition\" occurrences and re-insert the bookmark as a\n#    collapsed range right after the fourth (last) one, i.e. before\n#    the trailing \".\".\n$findEdition = $d.Content.Find\n$findEdition.Text = \"edition\"\n$findEdition.MatchCase = $true\n$findEdition.MatchWholeWord = $true\n\n$matchCount = 0\n$lastMatch = $null\nwhile ($findEdition.Execute()) {\n    $matchCount = $matchCount + 1\n    $lastMatch = $d.Range($findEdition.Parent.Start, $findEdition.Parent.End)\n}\n\n$insertionPoint = $d.Range($lastMatch.End, $lastMatch.End)\n$d.Bookmarks.Add(\"_GoBack\", $insertionPoint)\n"}
